$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Participantes Presentes",
    "Beatriz",
    "Daniela",
    "Gabriel",
    "João",
    "Kauê",
    "Luiz"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Columns.Item(1).ColumnWidth = 19

[void]$ws.Range("B8").Select()
